# Auto-generated Excel COM-interop edit script
# Updates market price snapshot cells (columns H:N) across the ALC, ARM, BSM, and GSM sheets
# to refresh currentAveragePrice / LevePrice / LeveProfit figures, per the scheduled price-data run.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(76, 8).Value = 2938.9683
$ws.Cells.Item(76, 9).Value = 2815.8164
$ws.Cells.Item(76, 10).Value = 3370
$ws.Cells.Item(76, 11).Value = 2815.8164
$ws.Cells.Item(76, 12).Value = 3370
$ws.Cells.Item(76, 13).Value = -2500.8164
$ws.Cells.Item(76, 14).Value = -4000

$ws.Cells.Item(79, 8).Value = 2938.9683
$ws.Cells.Item(79, 9).Value = 2815.8164
$ws.Cells.Item(79, 10).Value = 3370
$ws.Cells.Item(79, 11).Value = 2815.8164
$ws.Cells.Item(79, 12).Value = 3370
$ws.Cells.Item(79, 13).Value = -1723.8164
$ws.Cells.Item(79, 14).Value = -5554


$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(32, 8).Value = 11846.888
$ws.Cells.Item(32, 9).Value = 6576.1387
$ws.Cells.Item(32, 10).Value = 68946.664
$ws.Cells.Item(32, 11).Value = 6576.1387
$ws.Cells.Item(32, 12).Value = 68946.664
$ws.Cells.Item(32, 13).Value = -6289.1387
$ws.Cells.Item(32, 14).Value = -69520.664

$ws.Cells.Item(88, 8).Value = 1781.6364
$ws.Cells.Item(88, 9).Value = 1749.8
$ws.Cells.Item(88, 10).Value = 2100
$ws.Cells.Item(88, 11).Value = 1749.8
$ws.Cells.Item(88, 12).Value = 2100
$ws.Cells.Item(88, 13).Value = -1343.8
$ws.Cells.Item(88, 14).Value = -2912

$ws.Cells.Item(91, 8).Value = 1781.6364
$ws.Cells.Item(91, 9).Value = 1749.8
$ws.Cells.Item(91, 10).Value = 2100
$ws.Cells.Item(91, 11).Value = 1749.8
$ws.Cells.Item(91, 12).Value = 2100
$ws.Cells.Item(91, 13).Value = -345.8
$ws.Cells.Item(91, 14).Value = -4908


$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(86, 8).Value = 1801.2
$ws.Cells.Item(86, 9).Value = 1801.2
$ws.Cells.Item(86, 10).Value = 0
$ws.Cells.Item(86, 11).Value = 1801.2
$ws.Cells.Item(86, 12).Value = 0
$ws.Cells.Item(86, 13).Value = -678.2
$ws.Cells.Item(86, 14).Value = $null

$ws.Cells.Item(89, 8).Value = 1801.2
$ws.Cells.Item(89, 9).Value = 1801.2
$ws.Cells.Item(89, 10).Value = 0
$ws.Cells.Item(89, 11).Value = 9006
$ws.Cells.Item(89, 12).Value = 0
$ws.Cells.Item(89, 13).Value = -3390
$ws.Cells.Item(89, 14).Value = $null

$ws.Cells.Item(117, 8).Value = 0
$ws.Cells.Item(117, 9).Value = 0
$ws.Cells.Item(117, 10).Value = 0
$ws.Cells.Item(117, 11).Value = 0
$ws.Cells.Item(117, 12).Value = 0

$ws.Cells.Item(118, 8).Value = 26867.334
$ws.Cells.Item(118, 9).Value = 0
$ws.Cells.Item(118, 10).Value = 26867.334
$ws.Cells.Item(118, 11).Value = 0
$ws.Cells.Item(118, 12).Value = 26867.334
$ws.Cells.Item(118, 14).Value = -30181.334

$ws.Cells.Item(119, 8).Value = 40500
$ws.Cells.Item(119, 9).Value = 0
$ws.Cells.Item(119, 10).Value = 40500
$ws.Cells.Item(119, 11).Value = 0
$ws.Cells.Item(119, 12).Value = 40500
$ws.Cells.Item(119, 14).Value = -50176

$ws.Cells.Item(120, 8).Value = 39000
$ws.Cells.Item(120, 9).Value = 0
$ws.Cells.Item(120, 10).Value = 39000
$ws.Cells.Item(120, 11).Value = 0
$ws.Cells.Item(120, 12).Value = 39000
$ws.Cells.Item(120, 14).Value = -48676

$ws.Cells.Item(122, 8).Value = 0
$ws.Cells.Item(122, 9).Value = 0
$ws.Cells.Item(122, 10).Value = 0
$ws.Cells.Item(122, 11).Value = 0
$ws.Cells.Item(122, 12).Value = 0

$ws.Cells.Item(123, 8).Value = 41140
$ws.Cells.Item(123, 9).Value = 0
$ws.Cells.Item(123, 10).Value = 41140
$ws.Cells.Item(123, 11).Value = 0
$ws.Cells.Item(123, 12).Value = 41140
$ws.Cells.Item(123, 14).Value = -50940

$ws.Cells.Item(124, 8).Value = 33990
$ws.Cells.Item(124, 9).Value = 0
$ws.Cells.Item(124, 10).Value = 33990
$ws.Cells.Item(124, 11).Value = 0
$ws.Cells.Item(124, 12).Value = 33990
$ws.Cells.Item(124, 14).Value = -43810

$ws.Cells.Item(125, 8).Value = 37500
$ws.Cells.Item(125, 9).Value = 0
$ws.Cells.Item(125, 10).Value = 37500
$ws.Cells.Item(125, 11).Value = 0
$ws.Cells.Item(125, 12).Value = 37500
$ws.Cells.Item(125, 14).Value = -47340

$ws.Cells.Item(126, 8).Value = 46765
$ws.Cells.Item(126, 9).Value = 0
$ws.Cells.Item(126, 10).Value = 46765
$ws.Cells.Item(126, 11).Value = 0
$ws.Cells.Item(126, 12).Value = 46765
$ws.Cells.Item(126, 14).Value = -56645

$ws.Cells.Item(127, 8).Value = 38497.5
$ws.Cells.Item(127, 9).Value = 0
$ws.Cells.Item(127, 10).Value = 38497.5
$ws.Cells.Item(127, 11).Value = 0
$ws.Cells.Item(127, 12).Value = 38497.5
$ws.Cells.Item(127, 14).Value = -48417.5

$ws.Cells.Item(128, 8).Value = 2680
$ws.Cells.Item(128, 9).Value = 2680
$ws.Cells.Item(128, 10).Value = 0
$ws.Cells.Item(128, 11).Value = 8040
$ws.Cells.Item(128, 12).Value = 0
$ws.Cells.Item(128, 13).Value = -5550

$ws.Cells.Item(129, 8).Value = 50709
$ws.Cells.Item(129, 9).Value = 0
$ws.Cells.Item(129, 10).Value = 50709
$ws.Cells.Item(129, 11).Value = 0
$ws.Cells.Item(129, 12).Value = 50709
$ws.Cells.Item(129, 14).Value = -60709

$ws.Cells.Item(130, 8).Value = 0
$ws.Cells.Item(130, 9).Value = 0
$ws.Cells.Item(130, 10).Value = 0
$ws.Cells.Item(130, 11).Value = 0
$ws.Cells.Item(130, 12).Value = 0

$ws.Cells.Item(131, 8).Value = 36197.5
$ws.Cells.Item(131, 9).Value = 0
$ws.Cells.Item(131, 10).Value = 36197.5
$ws.Cells.Item(131, 11).Value = 0
$ws.Cells.Item(131, 12).Value = 36197.5
$ws.Cells.Item(131, 14).Value = -46277.5

$ws.Cells.Item(132, 8).Value = 0
$ws.Cells.Item(132, 9).Value = 0
$ws.Cells.Item(132, 10).Value = 0
$ws.Cells.Item(132, 11).Value = 0
$ws.Cells.Item(132, 12).Value = 0

$ws.Cells.Item(133, 8).Value = 28000
$ws.Cells.Item(133, 9).Value = 0
$ws.Cells.Item(133, 10).Value = 28000
$ws.Cells.Item(133, 11).Value = 0
$ws.Cells.Item(133, 12).Value = 28000
$ws.Cells.Item(133, 14).Value = -38120

$ws.Cells.Item(134, 8).Value = 4775.9414
$ws.Cells.Item(134, 9).Value = 1950.4839
$ws.Cells.Item(134, 10).Value = 9155.4
$ws.Cells.Item(134, 11).Value = 5851.4517
$ws.Cells.Item(134, 12).Value = 27466.2
$ws.Cells.Item(134, 13).Value = -3316.4517
$ws.Cells.Item(134, 14).Value = -32536.2

$ws.Cells.Item(135, 8).Value = 0
$ws.Cells.Item(135, 9).Value = 0
$ws.Cells.Item(135, 10).Value = 0
$ws.Cells.Item(135, 11).Value = 0
$ws.Cells.Item(135, 12).Value = 0

$ws.Cells.Item(137, 8).Value = 50487.5
$ws.Cells.Item(137, 9).Value = 0
$ws.Cells.Item(137, 10).Value = 50487.5
$ws.Cells.Item(137, 11).Value = 0
$ws.Cells.Item(137, 12).Value = 50487.5
$ws.Cells.Item(137, 14).Value = -60687.5

$ws.Cells.Item(138, 8).Value = 0
$ws.Cells.Item(138, 9).Value = 0
$ws.Cells.Item(138, 10).Value = 0
$ws.Cells.Item(138, 11).Value = 0
$ws.Cells.Item(138, 12).Value = 0

$ws.Cells.Item(139, 8).Value = 0
$ws.Cells.Item(139, 9).Value = 0
$ws.Cells.Item(139, 10).Value = 0
$ws.Cells.Item(139, 11).Value = 0
$ws.Cells.Item(139, 12).Value = 0

$ws.Cells.Item(140, 8).Value = 22000
$ws.Cells.Item(140, 9).Value = 0
$ws.Cells.Item(140, 10).Value = 22000
$ws.Cells.Item(140, 11).Value = 0
$ws.Cells.Item(140, 12).Value = 22000
$ws.Cells.Item(140, 14).Value = -32360

$ws.Cells.Item(141, 8).Value = 50000
$ws.Cells.Item(141, 9).Value = 0
$ws.Cells.Item(141, 10).Value = 50000
$ws.Cells.Item(141, 11).Value = 0
$ws.Cells.Item(141, 12).Value = 50000
$ws.Cells.Item(141, 14).Value = -60360


$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(70, 8).Value = 6589.7
$ws.Cells.Item(70, 9).Value = 5055.4287
$ws.Cells.Item(70, 10).Value = 10169.667
$ws.Cells.Item(70, 11).Value = 5055.4287
$ws.Cells.Item(70, 12).Value = 10169.667
$ws.Cells.Item(70, 13).Value = -4785.4287
$ws.Cells.Item(70, 14).Value = -10709.667

$ws.Cells.Item(73, 8).Value = 6589.7
$ws.Cells.Item(73, 9).Value = 5055.4287
$ws.Cells.Item(73, 10).Value = 10169.667
$ws.Cells.Item(73, 11).Value = 5055.4287
$ws.Cells.Item(73, 12).Value = 10169.667
$ws.Cells.Item(73, 13).Value = -4119.4287
$ws.Cells.Item(73, 14).Value = -12041.667

